# Add a new "Image Url" column (V) with company-logo file names, matching
# the "show proper company logo" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column V: widen it like the other data columns ---
$ws.Columns.Item(22).ColumnWidth = 21.77734375

# --- Copy the existing bordered / wrap-text cell format (used by A4:U6)
#     onto the new column's cells so the new cells share that style. ---
$ws.Range("A4").Copy()
$ws.Range("V1:V6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Header ---
$ws.Range("V1").Value = "Image Url"

# --- Per-company logo file names (row order matches existing data rows) ---
$ws.Range("V2").Value = "Byjus.png"
$ws.Range("V3").Value = "Swiggy Access.png"
$ws.Range("V4").Value = "Lenskart.png"
$ws.Range("V5").Value = "Mamaearth.png"
$ws.Range("V6").Value = "zomato.png"

# --- Row heights: rows 1 & 2 grow slightly, row 3 grows a lot more to fit
#     the new wrapped cell content (row 3's bottom thick border is retained
#     automatically since we are not touching its borders). ---
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 43.8

# --- Reselect the new column, mirroring what the author did after adding it ---
$ws.Range("V1:V6").Select()
